$wb = $excel.ActiveWorkbook

# --- Trips sheet ---
$tripsWs = $wb.Worksheets.Item("Trips")

# TripId 90: remove "pera" traveller from the Travellers JSON list
$tripsWs.Range("G2").Value = '[{"TravellerId":108,"Email":"mika@gmail.com","FirstName":"mika"}]'

# TripId 91: Destination country corrected from "novann" to "Nova drzava"
$tripsWs.Range("F3").Value = '{"DestinationId":112,"City":"nova Destinacija","Country":"Nova drzava"}'
# TripId 91: Reviews cleared out
$tripsWs.Range("H3").Value = '[]'

# TripId 92: Destination changed from Paris to London
$tripsWs.Range("F4").Value = '{"DestinationId":113,"City":"London","Country":"England"}'

# --- Destinations sheet ---
$destWs = $wb.Worksheets.Item("Destinations")

# DestinationId 110 (Paris): Attractions cleared out
$destWs.Range("F2").Value = '[]'

# DestinationId 112: Country corrected from "novann" to "Nova drzava"
$destWs.Range("C4").Value = 'Nova drzava'

# --- View / selection changes ---
# Destinations becomes non-active tab with C4 selected
$destWs.Range("C4").Select()

# Trips becomes the active tab with F3 selected
$tripsWs.Range("F3").Select()
